# The sheet had an extra leading column A (values 0,4,11,14, header-bordered
# style) that needs to be removed entirely, shifting columns B:F left to
# become A:E. Deleting the whole column reproduces exactly that shift,
# including re-numbering the header row and all data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).Delete()
